# The source workbook had a "label" column (D) identifying each data row
# with a pcap file name (dns.pcap, SkypeIRC.cap, sip_rtp.pcap, netflow.pcap).
# Those rows/labels are no longer used anywhere (the chart reads its
# categories from column E), so this cleans them out, which in turn lets
# the shared-string table shrink to only the strings that are still
# referenced (Count, Time (s), Improved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused row labels in column D (rows 5-8).
$ws.Range("D5:D8").ClearContents()

# Move the selection the way the author left it.
$ws.Range("O12").Select()
